# Insert a new weekly price record for "Ajo" (Chino / Primera) at row 137.
# This shifts all existing rows from 137-159 down to 138-160 (matching the
# source diff, which effectively inserts one new data row and pushes the
# remaining rows down by one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 137, shifting rows 137..159 down to 138..160.
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A137").Value = 8
$ws.Range("B137").Value = "Terminal La Palmera de La Serena"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44491
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 100112003
$ws.Range("G137").Value = "Ajo"
$ws.Range("H137").Value = "Chino"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 700
$ws.Range("K137").Value = 19000
$ws.Range("L137").Value = 20000
$ws.Range("M137").Value = 19500
$ws.Range("N137").Value = "`$/caja 10 kilos"
$ws.Range("O137").Value = "China"
$ws.Range("P137").Value = 1950
$ws.Range("Q137").Value = 10
$ws.Range("R137").Value = "Hortaliza"
